$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46056.01041666666, 267, 971, 1465.585359184455, 1238),
    @(3, 46056.02083333334, 256, 955, 1460.395036285225, 1211),
    @(4, 46056.03125, 255, 0, 1455.204713385996, 255),
    @(5, 46056.04166666666, 253, 954, 1450.014390486765, 1207),
    @(6, 46056.05208333334, 258, 938, 1440.995464384301, 1196),
    @(7, 46056.0625, 255, 936, 1431.976538281837, 1191),
    @(8, 46056.07291666666, 256, 937, 1422.957612179373, 1193),
    @(9, 46056.08333333334, 255, 935, 1413.938686076909, 1190),
    @(10, 46056.09375, 256, 880, 1412.547363778826, 1136),
    @(11, 46056.10416666666, 255, 877, 1411.156041480744, 1132),
    @(12, 46056.11458333334, 256, 879, 1409.764719182661, 1135),
    @(13, 46056.125, 0, 877, 1408.373396884579, 877),
    @(14, 46056.13541666666, 0, 883, 1425.118602131714, 883),
    @(15, 46056.14583333334, 255, 882, 1441.863807378849, 1137),
    @(16, 46056.15625, 256, 883, 1458.609012625985, 1139),
    @(17, 46056.16666666666, 254, 884, 1475.35421787312, 1138),
    @(18, 46056.17708333334, 257, 901, 1497.191070679196, 1158),
    @(19, 46056.1875, 256, 0, 1519.027923485271, 256),
    @(20, 46056.19791666666, 258, 0, 1540.864776291347, 258),
    @(21, 46056.20833333334, 268, 913, 1562.701629097423, 1181),
    @(22, 46056.21875, 398, 1235, 1700.472139209739, 1633),
    @(23, 46056.22916666666, 404, 1242, 1838.242649322056, 1646),
    @(24, 46056.23958333334, 410, 1243, 1976.013159434372, 1653),
    @(25, 46056.25, 439, 1252, 2113.783669546688, 1691),
    @(26, 46056.26041666666, 808, 1439, 2191.076063595363, 2247),
    @(27, 46056.27083333334, 803, 1441, 2268.368457644038, 2244),
    @(28, 46056.28125, 834, 1444, 2345.660851692714, 2278),
    @(29, 46056.29166666666, 897, 1467, 2422.953245741388, 2364),
    @(30, 46056.30208333334, 830, 1601, 2424.916885013216, 2431),
    @(31, 46056.3125, 944, 1603, 2426.880524285044, 2547),
    @(32, 46056.32291666666, 947, 0, 2428.844163556872, 947),
    @(33, 46056.33333333334, 902, 1608, 2430.8078028287, 2510),
    @(34, 46056.34375, 938, 1603, 2427.936350426274, 2541),
    @(35, 46056.35416666666, 873, 1601, 2425.064898023848, 2474),
    @(36, 46056.36458333334, 786, 1595, 2422.193445621423, 2381),
    @(37, 46056.375, 775, 1593, 2419.321993218997, 2368),
    @(38, 46056.38541666666, 678, 1617, 2394.554482806322, 2295),
    @(39, 46056.39583333334, 675, 1620, 2369.786972393646, 2295),
    @(40, 46056.40625, 673, 1621, 2345.019461980972, 2294),
    @(41, 46056.41666666666, 676, 1622, 2320.251951568297, 2298),
    @(42, 46056.42708333334, 615, 1617, 2292.168949721953, 2232),
    @(43, 46056.4375, 610, 1611, 2264.085947875609, 2221),
    @(44, 46056.44791666666, 595, 1605, 2236.002946029265, 2200),
    @(45, 46056.45833333334, 593, 1613, 2207.919944182922, 2206),
    @(46, 46056.46875, 585, 1592, 2192.960762594967, 2177),
    @(47, 46056.47916666666, 580, 1584, 2178.001581007013, 2164),
    @(48, 46056.48958333334, 579, 1581, 2163.042399419058, 2160),
    @(49, 46056.5, 0, 1582, 2148.083217831103, 1582),
    @(50, 46056.51041666666, 0, 0, 2147.71811907209, 0),
    @(51, 46056.52083333334, 0, 0, 2147.353020313078, 0),
    @(52, 46056.53125, 0, 0, 2146.987921554065, 0),
    @(53, 46056.54166666666, 0, 0, 2146.622822795052, 0),
    @(54, 46056.55208333334, 0, 0, 2158.375055823418, 0),
    @(55, 46056.5625, 0, 0, 2170.127288851784, 0),
    @(56, 46056.57291666666, 0, 0, 2181.879521880151, 0),
    @(57, 46056.58333333334, 0, 0, 2193.631754908517, 0),
    @(58, 46056.59375, 0, 0, 2231.75990639652, 0),
    @(59, 46056.60416666666, 0, 0, 2269.888057884523, 0),
    @(60, 46056.61458333334, 0, 0, 2308.016209372525, 0),
    @(61, 46056.625, 0, 0, 2346.144360860528, 0),
    @(62, 46056.63541666666, 0, 0, 2383.996221395329, 0),
    @(63, 46056.64583333334, 0, 0, 2421.848081930129, 0),
    @(64, 46056.65625, 0, 0, 2459.69994246493, 0),
    @(65, 46056.66666666666, 0, 0, 2497.551802999731, 0),
    @(66, 46056.67708333334, 0, 0, 2524.687521579559, 0),
    @(67, 46056.6875, 0, 0, 2551.823240159387, 0),
    @(68, 46056.69791666666, 0, 0, 2578.958958739215, 0),
    @(69, 46056.70833333334, 0, 0, 2606.094677319043, 0),
    @(70, 46056.71875, 0, 0, 2611.038311867691, 0),
    @(71, 46056.72916666666, 0, 0, 2615.98194641634, 0),
    @(72, 46056.73958333334, 0, 0, 2620.925580964989, 0),
    @(73, 46056.75, 0, 0, 2625.869215513638, 0),
    @(74, 46056.76041666666, 0, 0, 2616.978567353478, 0),
    @(75, 46056.77083333334, 0, 0, 2608.087919193318, 0),
    @(76, 46056.78125, 0, 0, 2599.197271033158, 0),
    @(77, 46056.79166666666, 0, 0, 2590.306622872998, 0),
    @(78, 46056.80208333334, 0, 0, 2574.113999531345, 0),
    @(79, 46056.8125, 0, 0, 2557.921376189692, 0),
    @(80, 46056.82291666666, 0, 0, 2541.728752848038, 0),
    @(81, 46056.83333333334, 0, 0, 2525.536129506385, 0),
    @(82, 46056.84375, 0, 0, 2488.177645672769, 0),
    @(83, 46056.85416666666, 0, 0, 2450.819161839152, 0),
    @(84, 46056.86458333334, 0, 0, 2413.460678005536, 0),
    @(85, 46056.875, 0, 0, 2376.102194171919, 0),
    @(86, 46056.88541666666, 0, 0, 2247.666371251166, 0),
    @(87, 46056.89583333334, 0, 0, 2119.230548330414, 0),
    @(88, 46056.90625, 0, 0, 1990.79472540966, 0),
    @(89, 46056.91666666666, 0, 0, 1862.358902488907, 0),
    @(90, 46056.92708333334, 0, 0, 1777.182348750924, 0),
    @(91, 46056.9375, 0, 0, 1692.00579501294, 0),
    @(92, 46056.94791666666, 0, 0, 1606.829241274957, 0),
    @(93, 46056.95833333334, 0, 0, 1521.652687536973, 0),
    @(94, 46056.96875, 0, 0, 1483.456518462315, 0),
    @(95, 46056.97916666666, 0, 0, 1445.260349387657, 0),
    @(96, 46056.98958333334, 0, 0, 1407.064180313, 0),
    @(97, 46057.0, 0, 0, 1368.868011238342, 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
}

Write-Host "Done"